$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added for "Feria Lagunitas de Puerto Montt - Repollo".
# It belongs chronologically/row-wise at row 229, pushing every existing row at/after
# 229 down by one (dimension grows from A1:R303 to A1:R304).
$ws.Rows(229).Insert()

$ws.Range("A229").Value2 = 4
$ws.Range("B229").Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Range("C229").Value2 = "Los Lagos"
$ws.Range("D229").Value2 = 44524
$ws.Range("E229").Value2 = 10
$ws.Range("F229").Value2 = 100112006
$ws.Range("G229").Value2 = "Repollo"
$ws.Range("H229").Value2 = "Crespo record"
$ws.Range("I229").Value2 = "Segunda"
$ws.Range("J229").Value2 = 100
$ws.Range("K229").Value2 = 1000
$ws.Range("L229").Value2 = 1000
$ws.Range("M229").Value2 = 1000
$ws.Range("N229").Value2 = '$/unidad'
$ws.Range("O229").Value2 = "Región Metropolitana"
$ws.Range("P229").Value2 = 1000
$ws.Range("Q229").Value2 = 1
$ws.Range("R229").Value2 = "Hortaliza"
